# Weekly fruit/vegetable price update:
# Insert a new weekly record as row 28 (pushing the existing rows 28-30 down
# to 29-31) in the Esparragos / Feria Lagunitas de Puerto Montt sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 28; this shifts the old
# rows 28, 29, 30 down to 29, 30, 31 (and the workbook's used range grows
# from R30 to R31 automatically).
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with this week's data point.
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = "11/16/2021"
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 300000000
$ws.Range("G28").Value = "Espárragos"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 360
$ws.Range("K28").Value = 1600
$ws.Range("L28").Value = 1600
$ws.Range("M28").Value = 1600
$ws.Range("N28").Value = "`$/kilo"
$ws.Range("O28").Value = "Provincia de Linares"
$ws.Range("P28").Value = 1600
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"
